$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the underlying input values; dependent formulas (C4, F4, C5, F5)
# recalculate automatically.
$ws.Range("F2").Value = 1791425
$ws.Range("C3").Value = 1801019
$ws.Range("F3").Value = 1744385

$excel.CalculateFullRebuild()
